# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") holds per-game strikeout counts.
# This run regenerates those values (rows 2-32) with freshly computed
# figures instead of the previous "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2 through 32, in order.
$kValues = @(0, 2, 1, 0, 1, 1, 2, 1, 1, 0, 2, 0, 0, 1, 1, 1, 2, 3, 1, 1, 1, 0, 1, 1, 0, 0, 2, 3, 0, 1, 0)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
